$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  B="2022-03-15 16:59:30.279215"; C=120},
    @{Row=3;  B="2022-03-15 16:59:35.102788"; C=119},
    @{Row=4;  B="2022-03-15 16:59:37.241446"; C=118},
    @{Row=5;  B="2022-03-15 16:59:38.902441"; C=117},
    @{Row=6;  B="2022-03-15 16:59:40.005529"; C=116},
    @{Row=7;  B="2022-03-15 16:59:41.003133"; C=115},
    @{Row=8;  B="2022-03-15 16:59:42.059789"; C=114},
    @{Row=9;  B="2022-03-15 16:59:43.057506"; C=113},
    @{Row=10; B="2022-03-15 16:59:44.279287"; C=112},
    @{Row=11; B="2022-03-15 16:59:45.165435"; C=111},
    @{Row=12; B="2022-03-15 16:59:46.874566"; C=110},
    @{Row=13; B="2022-03-15 16:59:50.840736"; C=109},
    @{Row=14; B="2022-03-15 16:59:52.000288"; C=108},
    @{Row=15; B="2022-03-15 16:59:53.028074"; C=107},
    @{Row=16; B="2022-03-15 16:59:55.712180"; C=106},
    @{Row=17; B="2022-03-15 16:59:57.396259"; C=105},
    @{Row=18; B="2022-03-15 16:59:58.823782"; C=104},
    @{Row=19; B="2022-03-15 17:00:00.180320"; C=103},
    @{Row=20; B="2022-03-15 17:00:01.755551"; C=102},
    @{Row=21; B="2022-03-15 17:00:03.309958"; C=101},
    @{Row=22; B="2022-03-15 17:00:04.700003"; C=100}
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
